$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, Coin (B), Link (C), Price (D), Volume1h (E)
$data = @(
    @(2, 'Bitcoin', 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc', '27.699.56', '  +0.22%  '),
    @(3, 'Ethereum', 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth', '1.853.00', '  +0.42%  '),
    @(4, 'TetherUSD', 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt', '1.003', '  +0.36%  '),
    @(5, 'BNB', 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb', '312.34', '  -0.63%  '),
    @(6, 'USDC', 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc', '1.002', '  +0.26%  '),
    @(7, 'XRP', 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp', '0.4266', '  +0.56%  '),
    @(8, 'Cardano', 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada', '0.3590', '  -1.63%  '),
    @(9, 'Dogecoin', 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge', '0.07293', '  +0.04%  '),
    @(10, 'Polygon', 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic', '0.8778', '  -1.48%  '),
    @(11, 'Solana', 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol', '20.77', '  +0.11%  '),
    @(12, 'WrappedEther', 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth', '1.827.58', '  -1.55%  '),
    @(13, 'Chainlink', 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link', '6.546', '  -0.33%  '),
    @(14, 'Polkadot', 'https://coinranking.com/coin/25W7FG7om+polkadot-dot', '5.336', '  -0.06%  '),
    @(15, 'TRON', 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx', '0.07004', '  +1.55%  '),
    @(16, 'BinanceUSD', 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd', '1.003', '  +0.15%  '),
    @(17, 'Litecoin', 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc', '79.59', '  +0.59%  '),
    @(18, 'ShibaInu', 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib', '0.000008948', '  +0.60%  '),
    @(19, 'Dai', 'https://coinranking.com/coin/MoTuySvg7+dai-dai', '1.003', '  +0.36%  '),
    @(20, 'Avalanche', 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax', '15.25', '  -1.16%  '),
    @(21, 'WrappedBTC', 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc', '27.800.34', '  +0.58%  '),
    @(22, 'Uniswap', 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni', '4.993', '  -0.08%  '),
    @(23, 'Cosmos', 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom', '10.41', '  -1.80%  '),
    @(24, 'WrappedliquidstakedEther2.0', 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth', '2.066.47', '  +0.13%  '),
    @(25, 'Toncoin', 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton', '1.991', '  +4.45%  '),
    @(26, 'Monero', 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr', '154.94', '  +0.61%  '),
    @(27, 'EthereumClassic', 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc', '18.47', '  -2.70%  '),
    @(28, 'BitcoinCash', 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch', '120.15', '  -1.54%  '),
    @(29, 'InternetComputer(DFINITY)', 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp', '5.246', '  -0.97%  '),
    @(30, 'LidoDAOToken', 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo', '1.876', '  -2.38%  '),
    @(31, 'Stellar', 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm', '0.08918', '  -0.08%  '),
    @(32, 'ImmutableX', 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx', '0.7578', '  -2.13%  '),
    @(33, 'HuobiToken', 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht', '2.961', '  +1.59%  '),
    @(34, 'Filecoin', 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil', '4.512', '  -1.38%  '),
    @(35, 'ARBITRUM', 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb', '1.122', '  +1.62%  '),
    @(36, 'Frax', 'https://coinranking.com/coin/KfWtaeV1W+frax-frax', '1.002', '  +0.31%  '),
    @(37, 'Hedera', 'https://coinranking.com/coin/jad286TjB+hedera-hbar', '0.05421', '  +0.57%  '),
    @(38, 'TrustWalletToken', 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt', '1.103', '  -0.13%  '),
    @(39, 'VeChain', 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet', '0.01927', '  -1.16%  '),
    @(40, 'MXToken', 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx', '2.822', '  +0.76%  '),
    @(41, 'Algorand', 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo', '0.1671', '  +0.65%  '),
    @(42, 'TheSandbox', 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand', '0.5084', '  -0.38%  '),
    @(43, 'FraxShare', 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs', '6.614', '  -4.40%  '),
    @(44, 'Aptos', 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt', '8.404', '  +1.29%  '),
    @(45, 'Cronos', 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro', '0.06525', '  -0.94%  '),
    @(46, 'Quant', 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt', '105.93', '  +1.31%  '),
    @(47, 'EnergySwap', 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens', '10.32', '  -1.29%  '),
    @(48, 'Decentraland', 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana', '0.4670', '  -1.51%  '),
    @(49, 'PaxDollar', 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp', '1.002', '  +0.37%  '),
    @(50, 'NEARProtocol', 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near', '1.620', '  -1.02%  '),
    @(51, 'RenderToken', 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr', '1.795', '  +2.15%  ')
)

# Column D (Price) values must stay as text (they were stored as inline
# strings in the original file), not get coerced into numbers. Force the
# whole column to a text number format once, up front, before assigning
# any values, so every cell in the range gets a single consistent style.
$ws.Range("D2:D51").NumberFormat = "@"

foreach ($entry in $data) {
    $row = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
    $ws.Cells.Item($row, 5).Value = $entry[4]
}
